$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (26-33), matching columns A-I:
# A = index (numeric), B = label, C = AAPC, D = IC_Inf, E = IC_Sup, F = 2019, G = 2020, H = 2021, I = delta
$newRows = @(
    @{ Row = 26; A = 24; B = "Admissions_eld";              C = -23.3918064291519;   D = -23.64087211592295; E = -23.14192834813216; F = 246136; G = 169402; H = 147314; I = -40.14934832775376 },
    @{ Row = 27; A = 25; B = "Admissions_non_eld";          C = -16.25632106201723;  D = -16.64355958483086; E = -15.86728359394149; F = 107455; G = 89628;  H = 75396;  I = -29.83481457354241 },
    @{ Row = 28; A = 26; B = "Admissions_uti_eld";          C = -4.255397002338013;  D = -5.312155718430278; E = -3.186844386062371;  F = 15747;  G = 16686;  H = 14390;  I = -8.617514447196291 },
    @{ Row = 29; A = 27; B = "Admissions_uti_non_eld";      C = 4.402436423497003;   D = 2.747967374614246;  E = 6.08354607563113;    F = 6855;   G = 8244;   H = 7504;   I = 9.467541940189642 },
    @{ Row = 30; A = 28; B = "Admissions_non_uti_eld";      C = -24.96317288138551;  D = -25.2186574591482;  E = -24.70681546064813;  F = 230389; G = 152716; H = 132924; I = -42.30453710897656 },
    @{ Row = 31; A = 29; B = "Admissions_non_uti_non_eld";  C = -17.93190964976411;  D = -18.32889565371637; E = -17.53299398061882;  F = 100600; G = 81384;  H = 67892;  I = -32.51292246520875 },
    @{ Row = 32; A = 30; B = "Deaths_eld";                  C = -11.50267438425282;  D = -12.08411280793379; E = -10.91739057322116;  F = 50437;  G = 43586;  H = 39580;  I = -21.525863949085 },
    @{ Row = 33; A = 31; B = "Deaths_non_eld";              C = -2.614763784128404;  D = -4.04817946983721;  E = -1.159934429385656;  F = 8900;   G = 8878;   H = 8437;   I = -5.202247191011236 }
)

# Copy the style of the last existing "index" cell (A25) so the new A-column
# cells (A26:A33) pick up the same bold/centered/bordered style (s="1").
$ws.Range("A25").Copy()
$ws.Range("A26:A33").PasteSpecial(-4122)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}

$ws.Range("A1:I33").Select()
